# ---------------------------------------------------------------------------
# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" positioned right after "总计" (i.e.
#    immediately before the existing "2022-Q2" sheet).
# 2. Populate it with the fund-holding breakdown for 2022-Q3, matching the
#    layout/style used by the other quarterly sheets.
# 3. Update the "总计" (summary) sheet: insert the new 2022-Q3 roll-up row
#    at the top of the data (row 2), shifting 2022-Q2 / 2022-Q1 / 2021-Q1
#    down by one row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Insert "2022-Q3" sheet right before "2022-Q2" -----------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$ws = $wb.Worksheets.Add($beforeSheet)
$ws.Name = "2022-Q3"

# --- 2. Fill in the 2022-Q3 fund-holding table ------------------------------
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "'001532"
$ws.Range("C2").Value = "华安文体健康主题灵活配置混合A"
$ws.Range("D2").Value = "'42.79"
$ws.Range("E2").Value = "'87.37"
$ws.Range("F2").Value = "'2.49"
$ws.Range("G2").Value = "'1.0655"
$ws.Range("H2").Value = 6
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "'002229"
$ws.Range("C3").Value = "华夏经济转型股票"
$ws.Range("D3").Value = "'12.47"
$ws.Range("E3").Value = "'87.35"
$ws.Range("F3").Value = "'3.47"
$ws.Range("G3").Value = "'0.4327"
$ws.Range("H3").Value = 6
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "'006868"
$ws.Range("C4").Value = "华夏科技成长股票"
$ws.Range("D4").Value = "'5.37"
$ws.Range("E4").Value = "'88.33"
$ws.Range("F4").Value = "'4.32"
$ws.Range("G4").Value = "'0.2320"
$ws.Range("H4").Value = 3
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "'013116"
$ws.Range("C5").Value = "华安文体健康主题灵活配置混合C"
$ws.Range("D5").Value = "'4.16"
$ws.Range("E5").Value = "'87.37"
$ws.Range("F5").Value = "'2.49"
$ws.Range("G5").Value = "'0.1036"
$ws.Range("H5").Value = 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "'013877"
$ws.Range("C6").Value = "财通资管新能源汽车混合C"
$ws.Range("D6").Value = "'0.95"
$ws.Range("E6").Value = "'94.48"
$ws.Range("F6").Value = "'5.26"
$ws.Range("G6").Value = "'0.0500"
$ws.Range("H6").Value = 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "'001352"
$ws.Range("C7").Value = "民生加银新战略灵活配置混合A"
$ws.Range("D7").Value = "'0.77"
$ws.Range("E7").Value = "'46.20"
$ws.Range("F7").Value = "'3.28"
$ws.Range("G7").Value = "'0.0253"
$ws.Range("H7").Value = 6
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "'013876"
$ws.Range("C8").Value = "财通资管新能源汽车混合A"
$ws.Range("D8").Value = "'0.13"
$ws.Range("E8").Value = "'94.48"
$ws.Range("F8").Value = "'5.26"
$ws.Range("G8").Value = "'0.0068"
$ws.Range("H8").Value = 7
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "'009054"
$ws.Range("C9").Value = "圆信永丰沣泰混合"
$ws.Range("D9").Value = "'0.23"
$ws.Range("E9").Value = "'26.81"
$ws.Range("F9").Value = "'1.06"
$ws.Range("G9").Value = "'0.0024"
$ws.Range("H9").Value = 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "'011391"
$ws.Range("C10").Value = "民生加银新战略灵活配置混合C"
$ws.Range("D10").Value = "'0.02"
$ws.Range("E10").Value = "'46.20"
$ws.Range("F10").Value = "'3.28"
$ws.Range("G10").Value = "'0.0007"
$ws.Range("H10").Value = 6

# Match header row (B1:H1) and index column (A2:A10) styling to the rest of
# the workbook by copying the formatting from the "总计" sheet's header style.
$styleSource = $wb.Worksheets.Item("总计").Range("B1")
$styleSource.Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$ws.Range("A2:A10").PasteSpecial(-4122)

# --- 3. Update the "总计" summary sheet -------------------------------------
$totalWs = $wb.Worksheets.Item("总计")

$totalWs.Range("B2").Value = "2022-Q3"
$totalWs.Range("C2").Value = 9
$totalWs.Range("D2").Value = 1.92

$totalWs.Range("B3").Value = "2022-Q2"
$totalWs.Range("C3").Value = 13
$totalWs.Range("D3").Value = 2.69

$totalWs.Range("B4").Value = "2022-Q1"
$totalWs.Range("C4").Value = 3
$totalWs.Range("D4").Value = 2.1

$totalWs.Range("A4").Copy()
$totalWs.Range("A5").PasteSpecial(-4122)
$totalWs.Range("A5").Value = 3
$totalWs.Range("B5").Value = "2021-Q1"
$totalWs.Range("C5").Value = 1
$totalWs.Range("D5").Value = 0.08
